$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE: 001 -> 004 (kept as text, "'" forces text so Excel
# doesn't coerce the numeric-looking string to a number)
$ws.Range("J2").Value = "'004"

# REPORT_DATE: 2018-12-31 -> 2019-09-30
$ws.Range("N2").Value = "2019-09-30 00:00:00"

# Updated financial figures for row 2
$ws.Range("O2").Value = 42055175.49
$ws.Range("P2").Value = 265462037.41
$ws.Range("Q2").Value = 215991348.45

# TOE_RATIO no longer reported -> blank out
$ws.Range("R2").ClearContents()

$ws.Range("S2").Value = 168275840.43
$ws.Range("T2").Value = 168275840.43

# OPERATE_EXPENSE_RATIO no longer reported -> blank out
$ws.Range("U2").ClearContents()

$ws.Range("V2").Value = 16473186.32
$ws.Range("W2").Value = 18564739.48
$ws.Range("X2").Value = 784201.02
$ws.Range("Y2").Value = 50272773.92
$ws.Range("Z2").Value = 49816459.45
$ws.Range("AA2").Value = 7761283.96

$ws.Range("AG2").Value = 3095204.65

# TOI_RATIO / OPERATE_PROFIT_RATIO / PARENT_NETPROFIT_RATIO no longer
# reported -> blank out
$ws.Range("AP2").ClearContents()
$ws.Range("AQ2").ClearContents()
$ws.Range("AR2").ClearContents()

$ws.Range("AS2").Value = 41543548

# DPN_RATIO no longer reported -> blank out
$ws.Range("AT2").ClearContents()
